$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -2
    3  = -6
    4  = 2
    5  = 6
    6  = 2
    7  = 2
    8  = 5
    9  = 4
    10 = 3
    11 = -2
    12 = 5
    13 = 6
    14 = -1
    15 = -2
    16 = 1
    17 = 4
    18 = 3
    19 = -1
    20 = 6
    21 = -5
    22 = 4
    23 = -4
    24 = -6
    25 = -3
    26 = 8
    27 = 1
    28 = 1
    29 = -2
    30 = -1
    31 = -3
    32 = -2
    34 = -1
    35 = -4
    36 = -2
    37 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
